# Update the "Presented By" student-name line on slide 1.
#   "Student Name-S Priya"  ->  "Student Name-JV Adlin Thirsha"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the textbox that holds the "Presented By / Student Name / College /
# Department" bullet list (falls back to the known shape if no match is
# found, e.g. "TextBox 3" in the original deck).
$sh = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $cand = $s.Shapes.Item($i)
    if ($cand.HasTextFrame -and $cand.TextFrame.TextRange.Text -match "Priya") {
        $sh = $cand
    }
}
if ($sh -eq $null) {
    $sh = $s.Shapes.Item(3)
}

$tr = $sh.TextFrame.TextRange

# In the full textbox string, paragraph 2 ("Student Name-S Priya") starts at
# character 15 (1-based): "Student Name-S " occupies chars 15-29 (15 chars)
# and "Priya" occupies chars 30-34 (5 chars).

# Edit right-to-left so earlier character offsets stay valid after each
# replacement changes the overall text length.

# 1) "Priya" -> "Adlin Thirsha"
$tr.Characters(30, 5).Text = "Adlin Thirsha"

# 2) "Student Name-S " -> "Student Name-JV "
$tr.Characters(15, 15).Text = "Student Name-JV "

# "Student Name-JV " is 16 characters, so the replaced "Adlin Thirsha" text
# now starts at character 31. Touch the character formatting at the word
# boundaries (re-asserting the existing font name is a no-op visually) so
# PowerPoint splits "Adlin Thirsha" into separate runs for "Adlin", " " and
# "Thirsha", matching how the words were actually typed/autocorrected.
$tr.Characters(36, 8).Font.Name = "Arial"   # " Thirsha"
$tr.Characters(37, 7).Font.Name = "Arial"   # "Thirsha"
